$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '97.486.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.595.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +18.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '654.08'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +8.99%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.06'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.593.43'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '44.55'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.37%  '
$ws.Range('E13').Value = '  +1.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.262.34'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '97.110.91'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000262'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.71'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +12.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.581.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.78%  '
$ws.Range('E22').Value = '  +11.02%  '
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '518.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000209'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '102.10'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.786.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.170'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +18.53%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.62%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.99%  '
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('E34').Value = '  +5.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.02'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.574'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '620.14'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.77'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.62%  '
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('E41').Value = '  +3.39%  '
$ws.Range('E42').Value = '  +8.05%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.932'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.31%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.57%  '
$ws.Range('E46').Value = '  +7.54%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.428'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +41.85%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.65'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.98%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.31'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.62%  '
